# Insert a new block of TBV-tag paragraphs right after the blank paragraph
# that precedes the existing "[PUMP:TBV:1]" paragraph (and before it).
#
# Target structure (matches the author's diff):
#   [PUMP:TBV:1111]                                          (Normal)
#   PUMP:HRD:3350                                             (List Bullet)
#   Details regarding the full color touchscreen.             (indent 360 twips / 18pt)
#   PUMP:HTP:1500                                              (indent 720 twips / 36pt)
#   Test 1500                                                  (indent 720 twips / 36pt)
#   PUMP:HRD:0000                                              (List Bullet)
#   Details regarding the size and weight of the pump.        (indent 360 twips / 18pt)

$d = $word.ActiveDocument

# Locate the "[PUMP:TBV:1]" paragraph by its text (robust to any surrounding
# content); the new block is inserted immediately before it, i.e. right after
# the blank paragraph that precedes it in the original document.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "[PUMP:TBV:1]`r") {
        $targetIndex = $i
        break
    }
}
$anchorIndex = $targetIndex - 1

# New paragraphs are inserted immediately after the anchor, each built by
# inserting a new paragraph mark after the previous paragraph and then
# filling in the freshly created paragraph's style/indent/text. (A newly
# inserted paragraph inherits the formatting of the one it follows, so the
# style is always set explicitly to avoid unwanted inheritance.)

$anchor = $d.Paragraphs($anchorIndex).Range
$anchor.InsertParagraphAfter()
$p1 = $d.Paragraphs($anchorIndex + 1)
$p1.Style = "Normal"
$p1.Range.Text = "[PUMP:TBV:1111]"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($anchorIndex + 2)
$p2.Style = "List Bullet"
$p2.Range.Text = "PUMP:HRD:3350"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($anchorIndex + 3)
$p3.Style = "Normal"
$p3.LeftIndent = 18
$p3.Range.Text = "Details regarding the full color touchscreen. "

$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($anchorIndex + 4)
$p4.Style = "Normal"
$p4.LeftIndent = 36
$p4.Range.Text = "PUMP:HTP:1500"

$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs($anchorIndex + 5)
$p5.Style = "Normal"
$p5.LeftIndent = 36
$p5.Range.Text = "Test 1500 "

$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs($anchorIndex + 6)
$p6.Style = "List Bullet"
$p6.Range.Text = "PUMP:HRD:0000"

$p6.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs($anchorIndex + 7)
$p7.Style = "Normal"
$p7.LeftIndent = 18
$p7.Range.Text = "Details regarding the size and weight of the pump. "
